$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.173589587211609
$ws.Range("B1").Value = 2.189854621887207
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.372869968414307
$ws.Range("E1").Value = 1.227529525756836
